$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.66779673140667
$ws.Range("C2").Value = 9.683012303539137
$ws.Range("D2").Value = 14.12482764326565
$ws.Range("E2").Value = 14.96331208326708
$ws.Range("G2").Value = 42.97036403438559
$ws.Range("H2").Value = 17.98307346718465
$ws.Range("J2").Value = 8.978851923871614
$ws.Range("K2").Value = 10.52403804654302
$ws.Range("L2").Value = 11.75773452428673
$ws.Range("M2").Value = 16.3195047908041
$ws.Range("N2").Value = 21.85342209167314
$ws.Range("O2").Value = 29.23482447139614
$ws.Range("B3").Value = 14.50943044968309
$ws.Range("C3").Value = 9.680791312146726
$ws.Range("D3").Value = 14.1277388665587
$ws.Range("E3").Value = 14.99185532894032
$ws.Range("G3").Value = 43.07098362960609
$ws.Range("H3").Value = 18.02914779633803
$ws.Range("J3").Value = 8.986292062133568
$ws.Range("K3").Value = 10.40321605761696
$ws.Range("L3").Value = 11.76785129926367
$ws.Range("M3").Value = 16.30277518744363
$ws.Range("N3").Value = 21.91238315001919
$ws.Range("O3").Value = 29.31224649744936
$ws.Range("B4").Value = 14.41392994885017
$ws.Range("C4").Value = 9.679666462962013
$ws.Range("D4").Value = 14.13185802867452
$ws.Range("E4").Value = 15.01103918380295
$ws.Range("G4").Value = 43.14166748079858
$ws.Range("H4").Value = 18.05967607799375
$ws.Range("J4").Value = 8.991104679527353
$ws.Range("K4").Value = 10.32995529304104
$ws.Range("L4").Value = 11.7753429235318
$ws.Range("M4").Value = 16.29458286514575
$ws.Range("N4").Value = 21.95029998693401
$ws.Range("O4").Value = 29.36434284883751
$ws.Range("B5").Value = 14.37549084897221
$ws.Range("C5").Value = 9.679268527821385
$ws.Range("D5").Value = 14.13412418352994
$ws.Range("E5").Value = 15.01927426599301
$ws.Range("G5").Value = 43.17270704389361
$ws.Range("H5").Value = 18.07267994382182
$ws.Range("J5").Value = 8.993127476250594
$ws.Range("K5").Value = 10.3003648499886
$ws.Range("L5").Value = 11.77871833338595
$ws.Range("M5").Value = 16.29177072522101
$ws.Range("N5").Value = 21.9661837677646
$ws.Range("O5").Value = 29.38671866934603
$ws.Range("B6").Value = 14.36913813500818
$ws.Range("C6").Value = 9.679206113377894
$ws.Range("D6").Value = 14.13453599651881
$ws.Range("E6").Value = 15.0206669270979
$ws.Range("G6").Value = 43.17799604056832
$ws.Range("H6").Value = 18.07487326403544
$ws.Range("J6").Value = 8.993467086692174
$ws.Range("K6").Value = 10.29546823023394
$ws.Range("L6").Value = 11.7792983155512
$ws.Range("M6").Value = 16.29133565134137
$ws.Range("N6").Value = 21.96884740740842
$ws.Range("O6").Value = 29.39050337317294
$ws.Range("B7").Value = 14.41340955578862
$ws.Range("C7").Value = 9.679660850993743
$ws.Range("D7").Value = 14.13188621031135
$ws.Range("E7").Value = 15.01114855382151
$ws.Range("G7").Value = 43.14207704526698
$ws.Range("H7").Value = 18.0598491711066
$ws.Range("J7").Value = 8.991131709957568
$ws.Range("K7").Value = 10.32955511615294
$ws.Range("L7").Value = 11.77538713878272
$ws.Range("M7").Value = 16.29454280445518
$ws.Range("N7").Value = 21.9505124489022
$ws.Range("O7").Value = 29.36463997667582
$ws.Range("B8").Value = 14.61285655387475
$ws.Range("C8").Value = 9.682197165925615
$ws.Range("D8").Value = 14.12534830540279
$ws.Range("E8").Value = 14.97280996581637
$ws.Range("G8").Value = 43.00320790206224
$ws.Range("H8").Value = 17.99849558493802
$ws.Range("J8").Value = 8.98136667533204
$ws.Range("K8").Value = 10.4822057470001
$ws.Range("L8").Value = 11.76095756343755
$ws.Range("M8").Value = 16.31330684556769
$ws.Range("N8").Value = 21.8733967650036
$ws.Range("O8").Value = 29.26057325087055
$ws.Range("B9").Value = 15.01590502864239
$ws.Range("C9").Value = 9.689049398043098
$ws.Range("D9").Value = 14.13096597176584
$ws.Range("E9").Value = 14.91076102295184
$ws.Range("G9").Value = 42.80166733947677
$ws.Range("H9").Value = 17.89592176331612
$ws.Range("J9").Value = 8.964148249057619
$ws.Range("K9").Value = 10.78749186121161
$ws.Range("L9").Value = 11.7427861331122
$ws.Range("M9").Value = 16.36644990182662
$ws.Range("N9").Value = 21.73571968405714
$ws.Range("O9").Value = 29.09268162676157
$ws.Range("B10").Value = 15.31679566280859
$ws.Range("C10").Value = 9.695207303832913
$ws.Range("D10").Value = 14.14623947589328
$ws.Range("E10").Value = 14.8731470488459
$ws.Range("G10").Value = 42.69691782277057
$ws.Range("H10").Value = 17.83134729911799
$ws.Range("J10").Value = 8.952663735670008
$ws.Range("K10").Value = 11.0135301357407
$ws.Range("L10").Value = 11.73556608594311
$ws.Range("M10").Value = 16.41523982077928
$ws.Range("N10").Value = 21.64274505530965
$ws.Range("O10").Value = 28.99140197263325
$ws.Range("B11").Value = 15.4541816115782
$ws.Range("C11").Value = 9.698248043417074
$ws.Range("D11").Value = 14.15558345530751
$ws.Range("E11").Value = 14.85775965538407
$ws.Range("G11").Value = 42.65870006441166
$ws.Range("H11").Value = 17.80430687443631
$ws.Range("J11").Value = 8.947689911683828
$ws.Range("K11").Value = 11.11634453853105
$ws.Range("L11").Value = 11.73360295338213
$ws.Range("M11").Value = 16.43950009498653
$ws.Range("N11").Value = 21.60220657346471
$ws.Range("O11").Value = 28.95012141214289
$ws.Range("B12").Value = 15.50623515102351
$ws.Range("C12").Value = 9.699433516425589
$ws.Range("D12").Value = 14.15946381966968
$ws.Range("E12").Value = 14.85218007288597
$ws.Range("G12").Value = 42.64558602998106
$ws.Range("H12").Value = 17.79440266379422
$ws.Range("J12").Value = 8.945842302615088
$ws.Range("K12").Value = 11.15524352909837
$ws.Range("L12").Value = 11.73304865881345
$ws.Range("M12").Value = 16.44897903131304
$ws.Range("N12").Value = 21.58710695337965
$ws.Range("O12").Value = 28.93517872789668
$ws.Range("B13").Value = 15.49502392951572
$ws.Range("C13").Value = 9.699176697779167
$ws.Range("D13").Value = 14.15861294602187
$ws.Range("E13").Value = 14.85337074617414
$ws.Range("G13").Value = 42.6483499436893
$ws.Range("H13").Value = 17.79652080058209
$ws.Range("J13").Value = 8.946238625437786
$ws.Range("K13").Value = 11.14686799104248
$ws.Range("L13").Value = 11.73315964058942
$ws.Range("M13").Value = 16.44692466251445
$ws.Range("N13").Value = 21.59034776214838
$ws.Range("O13").Value = 28.93836623874342
$ws.Range("B14").Value = 15.45846377115598
$ws.Range("C14").Value = 9.698344893817913
$ws.Range("D14").Value = 14.15589585634791
$ws.Range("E14").Value = 14.85729566682413
$ws.Range("G14").Value = 42.65759393652321
$ws.Range("H14").Value = 17.80348532883846
$ws.Range("J14").Value = 8.947537189777165
$ws.Range("K14").Value = 11.11954566138405
$ws.Range("L14").Value = 11.73355356669273
$ws.Range("M14").Value = 16.44027410773718
$ws.Range("N14").Value = 21.60095928604649
$ws.Range("O14").Value = 28.94887825119277
$ws.Range("B15").Value = 15.43607200802806
$ws.Range("C15").Value = 9.697839806398523
$ws.Range("D15").Value = 14.15427601996029
$ws.Range("E15").Value = 14.85973198239929
$ws.Range("G15").Value = 42.66343306910164
$ws.Range("H15").Value = 17.80779497402722
$ws.Range("J15").Value = 8.94833726464427
$ws.Range("K15").Value = 11.10280447476129
$ws.Range("L15").Value = 11.7338194568173
$ws.Range("M15").Value = 16.43623834053888
$ws.Range("N15").Value = 21.60749186085323
$ws.Range("O15").Value = 28.9554069455561
$ws.Range("B16").Value = 15.30782357393054
$ws.Range("C16").Value = 9.695013366955312
$ws.Range("D16").Value = 14.14567682661405
$ws.Range("E16").Value = 14.87418728508112
$ws.Range("G16").Value = 42.69960542221236
$ws.Range("H16").Value = 17.83316140860839
$ws.Range("J16").Value = 8.952993814787627
$ws.Range("K16").Value = 11.00680791031259
$ws.Range("L16").Value = 11.73572089335933
$ws.Range("M16").Value = 16.41369550281816
$ws.Range("N16").Value = 21.64542958356038
$ws.Range("O16").Value = 28.99419620781066
$ws.Range("B17").Value = 15.2292456997589
$ws.Range("C17").Value = 9.693340472251451
$ws.Range("D17").Value = 14.14101343789969
$ws.Range("E17").Value = 14.88349617294019
$ws.Range("G17").Value = 42.7242134532366
$ws.Range("H17").Value = 17.84932066204811
$ws.Range("J17").Value = 8.955914514259188
$ws.Range("K17").Value = 10.94789053695122
$ws.Range("L17").Value = 11.73722515469285
$ws.Range("M17").Value = 16.40039180107649
$ws.Range("N17").Value = 21.66915214622183
$ws.Range("O17").Value = 29.01921974096854
$ws.Range("B18").Value = 15.18410036656297
$ws.Range("C18").Value = 9.692400833925452
$ws.Range("D18").Value = 14.13855686824264
$ws.Range("E18").Value = 14.88901264556686
$ws.Range("G18").Value = 42.7392551793268
$ws.Range("H18").Value = 17.8588348218132
$ws.Range("J18").Value = 8.957618014537113
$ws.Range("K18").Value = 10.91400366955105
$ws.Range("L18").Value = 11.73821479018867
$ws.Range("M18").Value = 16.39293458749025
$ws.Range("N18").Value = 21.68296209887228
$ws.Range("O18").Value = 29.03406365421273
$ws.Range("B19").Value = 15.16882499461808
$ws.Range("C19").Value = 9.692086577244433
$ws.Range("D19").Value = 14.13776395163133
$ws.Range("E19").Value = 14.89090831039273
$ws.Range("G19").Value = 42.74450047653266
$ws.Range("H19").Value = 17.8620939111895
$ws.Range("J19").Value = 8.958198847079464
$ws.Range("K19").Value = 10.90253131810689
$ws.Range("L19").Value = 11.73857126040046
$ws.Range("M19").Value = 16.39044329306886
$ws.Range("N19").Value = 21.68766634658911
$ws.Range("O19").Value = 29.03916700517001
$ws.Range("B20").Value = 15.23760554001384
$ws.Range("C20").Value = 9.693516222277431
$ws.Range("D20").Value = 14.14148652335867
$ws.Range("E20").Value = 14.88248843784204
$ws.Range("G20").Value = 42.72150198460131
$ws.Range("H20").Value = 17.84757773745918
$ws.Range("J20").Value = 8.955601160357336
$ws.Range("K20").Value = 10.95416254487148
$ws.Range("L20").Value = 11.73705215244852
$ws.Range("M20").Value = 16.40178788453406
$ws.Range("N20").Value = 21.6666097326322
$ws.Range("O20").Value = 29.01650926012584
$ws.Range("B21").Value = 15.46920195681426
$ws.Range("C21").Value = 9.698588295170666
$ws.Range("D21").Value = 14.15668467098563
$ws.Range("E21").Value = 14.85613611543617
$ws.Range("G21").Value = 42.65484187942295
$ws.Range("H21").Value = 17.8014305777643
$ws.Range("J21").Value = 8.94715479769223
$ws.Range("K21").Value = 11.12757208980285
$ws.Range("L21").Value = 11.73343273630965
$ws.Range("M21").Value = 16.44221964982927
$ws.Range("N21").Value = 21.59783560947479
$ws.Range("O21").Value = 28.94577190717675
$ws.Range("B22").Value = 15.62070926920771
$ws.Range("C22").Value = 9.702101229334028
$ws.Range("D22").Value = 14.16860964833363
$ws.Range("E22").Value = 14.84035448435927
$ws.Range("G22").Value = 42.61919272488211
$ws.Range("H22").Value = 17.77322561924265
$ws.Range("J22").Value = 8.941843618146157
$ws.Range("K22").Value = 11.24068885690441
$ws.Range("L22").Value = 11.73216901338149
$ws.Range("M22").Value = 16.47034426369374
$ws.Range("N22").Value = 21.5543527244575
$ws.Range("O22").Value = 28.90355920679263
$ws.Range("B23").Value = 15.53984823909063
$ws.Range("C23").Value = 9.700208331132236
$ws.Range("D23").Value = 14.16206368370879
$ws.Range("E23").Value = 14.84864575349853
$ws.Range("G23").Value = 42.63749450129433
$ws.Range("H23").Value = 17.78810037553665
$ws.Range("J23").Value = 8.944659222004262
$ws.Range("K23").Value = 11.1803468062317
$ws.Range("L23").Value = 11.73274298342493
$ws.Range("M23").Value = 16.45517975272952
$ws.Range("N23").Value = 21.57742667609336
$ws.Range("O23").Value = 28.92572117431443
$ws.Range("B24").Value = 15.2338259589487
$ws.Range("C24").Value = 9.693436696726913
$ws.Range("D24").Value = 14.14127194202924
$ws.Range("E24").Value = 14.88294352191868
$ws.Range("G24").Value = 42.72272505420756
$ws.Range("H24").Value = 17.84836501592325
$ws.Range("J24").Value = 8.955742751793556
$ws.Range("K24").Value = 10.95132701268574
$ws.Range("L24").Value = 11.7371299779115
$ws.Range("M24").Value = 16.40115611907365
$ws.Range("N24").Value = 21.6677586234025
$ws.Range("O24").Value = 29.01773324329284
$ws.Range("B25").Value = 14.90585012012126
$ws.Range("C25").Value = 9.686996797060297
$ws.Range("D25").Value = 14.12748112220537
$ws.Range("E25").Value = 14.9261442780822
$ws.Range("G25").Value = 42.84859232124218
$ws.Range("H25").Value = 17.92177441568259
$ws.Range("J25").Value = 8.968600748621284
$ws.Range("K25").Value = 10.70446212390783
$ws.Range("L25").Value = 11.7466221982816
$ws.Range("M25").Value = 16.35034495385107
$ws.Range("N25").Value = 21.77152320746545
$ws.Range("O25").Value = 29.13422567145
